# Apply the table-style GUID change to the three tables that used the
# old custom "{D90F27A0-ABB0-48A9-92B1-B9AACE924EFE}" table style, and
# re-point them at the built-in style
# "{69B1A8C9-525D-4081-BACB-8E3B3438BC72}" (slides 14, 15 and 16 each
# contain exactly one table, as shape 1).
$p = $ppt.ActivePresentation

$newTableStyle = "{69B1A8C9-525D-4081-BACB-8E3B3438BC72}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $tableShape = $slide.Shapes.Item(1)
    $tableShape.Table.ApplyStyle($newTableStyle)
}

# Swap the presentation's colour theme from the "Integral" / "Red
# Violet" palette over to the plain "Office Theme" / "Office" palette
# (the theme used by the slide master - and therefore every slide -
# changes; this mirrors the source deck's theme1.xml / theme2.xml
# content swap for the part that actually drives on-slide appearance).
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
